# Jimenez 2019 pattern generator / evaluation data refresh.
# The generator was re-run, producing a different (but equivalent) traversal
# order of the 40 possible W/S note patterns; columns D-H (RC2, RC3,
# Complexitat, RC2 posicions, RC3 posicions) are recomputed per row to stay
# consistent with the new "Exemple" (B) sequence on each line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Exemple 2): ATAAATATAT
$ws.Cells.Item(3, 2).Value = "ATAAATATAT"
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "3"
$ws.Cells.Item(3, 8).Value = ""

# Row 4 (Exemple 3): AAATATATAT
$ws.Cells.Item(4, 2).Value = "AAATATATAT"
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "1"
$ws.Cells.Item(4, 8).Value = ""

# Row 6 (Exemple 5): ATATTAATAT
$ws.Cells.Item(6, 2).Value = "ATATTAATAT"
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = ""
$ws.Cells.Item(6, 8).NumberFormat = "@"
$ws.Cells.Item(6, 8).Value = "5"

# Row 7 (Exemple 6): ATATAAATAT
$ws.Cells.Item(7, 2).Value = "ATATAAATAT"
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = "5"
$ws.Cells.Item(7, 8).Value = ""

# Row 8 (Exemple 7): TAATATATAT
$ws.Cells.Item(8, 2).Value = "TAATATATAT"
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = ""
$ws.Cells.Item(8, 8).NumberFormat = "@"
$ws.Cells.Item(8, 8).Value = "1"

# Row 10 (Exemple 9): ATATATAAAT
$ws.Cells.Item(10, 2).Value = "ATATATAAAT"
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = "7"
$ws.Cells.Item(10, 8).Value = ""

# Row 11 (Exemple 10): TAAAATATAT
$ws.Cells.Item(11, 2).Value = "TAAAATATAT"
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 2
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = "3"
$ws.Cells.Item(11, 8).NumberFormat = "@"
$ws.Cells.Item(11, 8).Value = "1"

# Row 12 (Exemple 11): AAATTAATAT
$ws.Cells.Item(12, 2).Value = "AAATTAATAT"
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 2
$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = "1"
$ws.Cells.Item(12, 8).NumberFormat = "@"
$ws.Cells.Item(12, 8).Value = "5"

# Row 13 (Exemple 12): TAATAAATAT
$ws.Cells.Item(13, 2).Value = "TAATAAATAT"
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 2
$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = "5"
$ws.Cells.Item(13, 8).NumberFormat = "@"
$ws.Cells.Item(13, 8).Value = "1"

# Row 14 (Exemple 13): AAATATAAAT
$ws.Cells.Item(14, 2).Value = "AAATATAAAT"
$ws.Cells.Item(14, 4).Value = 2
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 2
$ws.Cells.Item(14, 7).Value = "1, 7"
$ws.Cells.Item(14, 8).Value = ""

# Row 15 (Exemple 14): ATAATAATAT
$ws.Cells.Item(15, 2).Value = "ATAATAATAT"
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 2
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "3"
$ws.Cells.Item(15, 8).NumberFormat = "@"
$ws.Cells.Item(15, 8).Value = "5"

# Row 16 (Exemple 15): TAATATAAAT
$ws.Cells.Item(16, 2).Value = "TAATATAAAT"
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 2
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "7"
$ws.Cells.Item(16, 8).NumberFormat = "@"
$ws.Cells.Item(16, 8).Value = "1"

# Row 17 (Exemple 16): ATATAATAAT
$ws.Cells.Item(17, 2).Value = "ATATAATAAT"
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 2
$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "5"
$ws.Cells.Item(17, 8).NumberFormat = "@"
$ws.Cells.Item(17, 8).Value = "7"

# Row 18 (Exemple 17): AAAAATATAT
$ws.Cells.Item(18, 2).Value = "AAAAATATAT"
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 2
$ws.Cells.Item(18, 7).Value = "1, 3"
$ws.Cells.Item(18, 8).Value = ""

# Row 19 (Exemple 18): AAATAAATAT
$ws.Cells.Item(19, 2).Value = "AAATAAATAT"
$ws.Cells.Item(19, 4).Value = 2
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 2
$ws.Cells.Item(19, 7).Value = "1, 5"
$ws.Cells.Item(19, 8).Value = ""

# Row 20 (Exemple 19): AAATATTAAT
$ws.Cells.Item(20, 2).Value = "AAATATTAAT"
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 6).Value = 2
$ws.Cells.Item(20, 7).NumberFormat = "@"
$ws.Cells.Item(20, 7).Value = "1"
$ws.Cells.Item(20, 8).NumberFormat = "@"
$ws.Cells.Item(20, 8).Value = "7"

# Row 21 (Exemple 20): ATTAAAATAT
$ws.Cells.Item(21, 2).Value = "ATTAAAATAT"
$ws.Cells.Item(21, 4).Value = 1
$ws.Cells.Item(21, 5).Value = 1
$ws.Cells.Item(21, 6).Value = 2
$ws.Cells.Item(21, 7).NumberFormat = "@"
$ws.Cells.Item(21, 7).Value = "5"
$ws.Cells.Item(21, 8).NumberFormat = "@"
$ws.Cells.Item(21, 8).Value = "3"

# Row 22 (Exemple 21): ATAAATAAAT
$ws.Cells.Item(22, 2).Value = "ATAAATAAAT"
$ws.Cells.Item(22, 4).Value = 2
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(22, 6).Value = 2
$ws.Cells.Item(22, 7).Value = "3, 7"
$ws.Cells.Item(22, 8).Value = ""

# Row 23 (Exemple 22): ATAAAAATAT
$ws.Cells.Item(23, 2).Value = "ATAAAAATAT"
$ws.Cells.Item(23, 4).Value = 2
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 2
$ws.Cells.Item(23, 7).Value = "3, 5"
$ws.Cells.Item(23, 8).Value = ""

# Row 24 (Exemple 23): ATATAAAAAT
$ws.Cells.Item(24, 2).Value = "ATATAAAAAT"
$ws.Cells.Item(24, 4).Value = 2
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 6).Value = 2
$ws.Cells.Item(24, 7).Value = "5, 7"
$ws.Cells.Item(24, 8).Value = ""

# Row 25 (Exemple 24): AATAATATAT
$ws.Cells.Item(25, 2).Value = "AATAATATAT"
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(25, 5).Value = 1
$ws.Cells.Item(25, 6).Value = 2
$ws.Cells.Item(25, 7).NumberFormat = "@"
$ws.Cells.Item(25, 7).Value = "1"
$ws.Cells.Item(25, 8).NumberFormat = "@"
$ws.Cells.Item(25, 8).Value = "3"

# Row 26 (Exemple 25): ATAAATTAAT
$ws.Cells.Item(26, 2).Value = "ATAAATTAAT"
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(26, 5).Value = 1
$ws.Cells.Item(26, 6).Value = 2
$ws.Cells.Item(26, 7).NumberFormat = "@"
$ws.Cells.Item(26, 7).Value = "3"
$ws.Cells.Item(26, 8).NumberFormat = "@"
$ws.Cells.Item(26, 8).Value = "7"

# Row 27 (Exemple 26): AAATAATAAT
$ws.Cells.Item(27, 2).Value = "AAATAATAAT"
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 5).Value = 1
$ws.Cells.Item(27, 6).Value = 3
$ws.Cells.Item(27, 7).Value = "1, 5"
$ws.Cells.Item(27, 8).NumberFormat = "@"
$ws.Cells.Item(27, 8).Value = "7"

# Row 28 (Exemple 27): AATAAAATAT
$ws.Cells.Item(28, 2).Value = "AATAAAATAT"
$ws.Cells.Item(28, 4).Value = 2
$ws.Cells.Item(28, 5).Value = 1
$ws.Cells.Item(28, 6).Value = 3
$ws.Cells.Item(28, 7).Value = "1, 5"
$ws.Cells.Item(28, 8).NumberFormat = "@"
$ws.Cells.Item(28, 8).Value = "3"

# Row 29 (Exemple 28): AAATAAAAAT
$ws.Cells.Item(29, 2).Value = "AAATAAAAAT"
$ws.Cells.Item(29, 4).Value = 3
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 3
$ws.Cells.Item(29, 7).Value = "1, 5, 7"
$ws.Cells.Item(29, 8).Value = ""

# Row 31 (Exemple 30): TAAAATAAAT
$ws.Cells.Item(31, 2).Value = "TAAAATAAAT"
$ws.Cells.Item(31, 4).Value = 2
$ws.Cells.Item(31, 5).Value = 1
$ws.Cells.Item(31, 6).Value = 3
$ws.Cells.Item(31, 7).Value = "3, 7"
$ws.Cells.Item(31, 8).NumberFormat = "@"
$ws.Cells.Item(31, 8).Value = "1"

# Row 32 (Exemple 31): ATAAAATAAT
$ws.Cells.Item(32, 2).Value = "ATAAAATAAT"
$ws.Cells.Item(32, 4).Value = 2
$ws.Cells.Item(32, 5).Value = 1
$ws.Cells.Item(32, 6).Value = 3
$ws.Cells.Item(32, 7).Value = "3, 5"
$ws.Cells.Item(32, 8).NumberFormat = "@"
$ws.Cells.Item(32, 8).Value = "7"

# Row 33 (Exemple 32): TAAAAAATAT
$ws.Cells.Item(33, 2).Value = "TAAAAAATAT"
$ws.Cells.Item(33, 4).Value = 2
$ws.Cells.Item(33, 5).Value = 1
$ws.Cells.Item(33, 6).Value = 3
$ws.Cells.Item(33, 7).Value = "3, 5"
$ws.Cells.Item(33, 8).NumberFormat = "@"
$ws.Cells.Item(33, 8).Value = "1"

# Row 34 (Exemple 33): ATAAAAAAAT
$ws.Cells.Item(34, 2).Value = "ATAAAAAAAT"
$ws.Cells.Item(34, 4).Value = 3
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 3
$ws.Cells.Item(34, 7).Value = "3, 5, 7"
$ws.Cells.Item(34, 8).Value = ""

# Row 35 (Exemple 34): AAAAATAAAT
$ws.Cells.Item(35, 2).Value = "AAAAATAAAT"
$ws.Cells.Item(35, 4).Value = 3
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 3
$ws.Cells.Item(35, 7).Value = "1, 3, 7"
$ws.Cells.Item(35, 8).Value = ""

# Row 36 (Exemple 35): TAATAAAAAT
$ws.Cells.Item(36, 2).Value = "TAATAAAAAT"
$ws.Cells.Item(36, 4).Value = 2
$ws.Cells.Item(36, 5).Value = 1
$ws.Cells.Item(36, 6).Value = 3
$ws.Cells.Item(36, 7).Value = "5, 7"
$ws.Cells.Item(36, 8).NumberFormat = "@"
$ws.Cells.Item(36, 8).Value = "1"

# Row 37 (Exemple 36): AAAAATTAAT
$ws.Cells.Item(37, 2).Value = "AAAAATTAAT"
$ws.Cells.Item(37, 4).Value = 2
$ws.Cells.Item(37, 5).Value = 1
$ws.Cells.Item(37, 6).Value = 3
$ws.Cells.Item(37, 7).Value = "1, 3"
$ws.Cells.Item(37, 8).NumberFormat = "@"
$ws.Cells.Item(37, 8).Value = "7"

# Row 38 (Exemple 37): AAAATAATAT
$ws.Cells.Item(38, 2).Value = "AAAATAATAT"
$ws.Cells.Item(38, 4).Value = 2
$ws.Cells.Item(38, 5).Value = 1
$ws.Cells.Item(38, 6).Value = 3
$ws.Cells.Item(38, 7).Value = "1, 3"
$ws.Cells.Item(38, 8).NumberFormat = "@"
$ws.Cells.Item(38, 8).Value = "5"
